$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 14012.5
$ws.Range("I7").Value = 8555
$ws.Range("J7").Value = 15831.667
$ws.Range("K7").Value = 8555
$ws.Range("L7").Value = 15831.667
$ws.Range("M7").Value = -8443
$ws.Range("N7").Value = -16055.667
$ws.Range("H14").Value = 14012.5
$ws.Range("I14").Value = 8555
$ws.Range("J14").Value = 15831.667
$ws.Range("K14").Value = 8555
$ws.Range("L14").Value = 15831.667
$ws.Range("M14").Value = -8364
$ws.Range("N14").Value = -16213.667
$ws.Range("H16").Value = 49990
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 49990
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 49990
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -50450
$ws.Range("H33").Value = 1689853
$ws.Range("I33").Value = 2252417.5
$ws.Range("K33").Value = 2252417.5
$ws.Range("M33").Value = -2252188.5
$ws.Range("H62").Value = 4561.2
$ws.Range("I62").Value = 4561.2
$ws.Range("K62").Value = 4561.2
$ws.Range("M62").Value = -3937.2
$ws.Range("H65").Value = 4561.2
$ws.Range("I65").Value = 4561.2
$ws.Range("K65").Value = 22806
$ws.Range("M65").Value = -19686
$ws.Range("H100").Value = 27922.975
$ws.Range("I100").Value = 37168.645
$ws.Range("K100").Value = 37168.645
$ws.Range("M100").Value = -36627.645
$ws.Range("H116").Value = 4864.2
$ws.Range("I116").Value = 3132.125
$ws.Range("K116").Value = 3132.125
$ws.Range("M116").Value = 309.875

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 1714.3
$ws.Range("I16").Value = 1723.8334
$ws.Range("J16").Value = 1700
$ws.Range("K16").Value = 1723.8334
$ws.Range("L16").Value = 1700
$ws.Range("M16").Value = -1436.8334
$ws.Range("N16").Value = -2274
$ws.Range("H61").Value = 7261.4116
$ws.Range("J61").Value = 21449
$ws.Range("L61").Value = 21449
$ws.Range("N61").Value = -21873
$ws.Range("H136").Value = 7261.4116
$ws.Range("J136").Value = 21449
$ws.Range("L136").Value = 64347
$ws.Range("N136").Value = -69447

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 372.39395
$ws.Range("I107").Value = 383.84616
$ws.Range("K107").Value = 383.84616
$ws.Range("M107").Value = 1536.15384

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 2981.8
$ws.Range("J8").Value = 2981.8
$ws.Range("L8").Value = 2981.8
$ws.Range("N8").Value = -3261.8
$ws.Range("H19").Value = 1878.2632
$ws.Range("I19").Value = 838
$ws.Range("J19").Value = 2485.0833
$ws.Range("K19").Value = 838
$ws.Range("L19").Value = 2485.0833
$ws.Range("M19").Value = -668
$ws.Range("N19").Value = -2825.0833
$ws.Range("H24").Value = 1878.2632
$ws.Range("I24").Value = 838
$ws.Range("J24").Value = 2485.0833
$ws.Range("K24").Value = 838
$ws.Range("L24").Value = 2485.0833
$ws.Range("M24").Value = -668
$ws.Range("N24").Value = -2825.0833
$ws.Range("H31").Value = 3353
$ws.Range("I31").Value = 2364.4546
$ws.Range("J31").Value = 4341.5454
$ws.Range("K31").Value = 2364.4546
$ws.Range("L31").Value = 4341.5454
$ws.Range("M31").Value = -2069.4546
$ws.Range("N31").Value = -4931.5454
$ws.Range("H34").Value = 3353
$ws.Range("I34").Value = 2364.4546
$ws.Range("J34").Value = 4341.5454
$ws.Range("K34").Value = 2364.4546
$ws.Range("L34").Value = 4341.5454
$ws.Range("M34").Value = -2162.4546
$ws.Range("N34").Value = -4745.5454
$ws.Range("H86").Value = 3583.842
$ws.Range("I86").Value = 3104.818
$ws.Range("K86").Value = 3104.818
$ws.Range("M86").Value = -1981.818
$ws.Range("H89").Value = 3583.842
$ws.Range("I89").Value = 3104.818
$ws.Range("K89").Value = 15524.09
$ws.Range("M89").Value = -9908.09
$ws.Range("H105").Value = 1329.7368
$ws.Range("I105").Value = 1355.5294
$ws.Range("K105").Value = 1355.5294
$ws.Range("M105").Value = 391.4706000000001
$ws.Range("H107").Value = 450.36365
$ws.Range("I107").Value = 289.33334
$ws.Range("K107").Value = 289.33334
$ws.Range("M107").Value = 1630.66666

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2678.2632
$ws.Range("I34").Value = 399.1111
$ws.Range("J34").Value = 4729.5
$ws.Range("K34").Value = 1197.3333
$ws.Range("L34").Value = 14188.5
$ws.Range("M34").Value = -1113.3333
$ws.Range("N34").Value = -14356.5
$ws.Range("H39").Value = 7142
$ws.Range("I39").Value = 1850
$ws.Range("J39").Value = 8024
$ws.Range("K39").Value = 5550
$ws.Range("L39").Value = 24072
$ws.Range("M39").Value = -5256
$ws.Range("N39").Value = -24660
$ws.Range("H55").Value = 5369
$ws.Range("J55").Value = 6390.6
$ws.Range("L55").Value = 19171.8
$ws.Range("N55").Value = -19525.8

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 32157
$ws.Range("J57").Value = 32157
$ws.Range("L57").Value = 32157
$ws.Range("N57").Value = -33797

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5271.75
$ws.Range("I7").Value = 5055.3
$ws.Range("K7").Value = 5055.3
$ws.Range("M7").Value = -4943.3
$ws.Range("H9").Value = 507.75
$ws.Range("I9").Value = 554.5714
$ws.Range("J9").Value = 180
$ws.Range("K9").Value = 554.5714
$ws.Range("L9").Value = 180
$ws.Range("M9").Value = -330.5714
$ws.Range("N9").Value = -628
$ws.Range("H18").Value = 900
$ws.Range("I18").Value = 900
$ws.Range("K18").Value = 900
$ws.Range("M18").Value = -728
$ws.Range("H19").Value = 6106
$ws.Range("J19").Value = 7466.3335
$ws.Range("L19").Value = 7466.3335
$ws.Range("N19").Value = -7806.3335
$ws.Range("H61").Value = 3568.5
$ws.Range("I61").Value = 3558.5386
$ws.Range("K61").Value = 3558.5386
$ws.Range("M61").Value = -3356.5386
$ws.Range("H68").Value = 3138.4167
$ws.Range("I68").Value = 2635.1614
$ws.Range("J68").Value = 6258.6
$ws.Range("K68").Value = 2635.1614
$ws.Range("L68").Value = 6258.6
$ws.Range("M68").Value = -1886.1614
$ws.Range("N68").Value = -7756.6
$ws.Range("H71").Value = 3138.4167
$ws.Range("I71").Value = 2635.1614
$ws.Range("J71").Value = 6258.6
$ws.Range("K71").Value = 13175.807
$ws.Range("L71").Value = 31293
$ws.Range("M71").Value = -9431.807000000001
$ws.Range("N71").Value = -38781
$ws.Range("H113").Value = 3568.5
$ws.Range("I113").Value = 3558.5386
$ws.Range("K113").Value = 3558.5386
$ws.Range("M113").Value = -1388.5386
$ws.Range("H122").Value = 5733.1113
$ws.Range("I122").Value = 5733.1113
$ws.Range("K122").Value = 17199.3339
$ws.Range("M122").Value = -14749.3339
$ws.Range("H126").Value = 5271.75
$ws.Range("I126").Value = 5055.3
$ws.Range("K126").Value = 15165.9
$ws.Range("M126").Value = -12695.9
$ws.Range("H132").Value = 6112.8
$ws.Range("I132").Value = 2227
$ws.Range("J132").Value = 9998.6
$ws.Range("K132").Value = 6681
$ws.Range("L132").Value = 29995.8
$ws.Range("M132").Value = -4151
$ws.Range("N132").Value = -35055.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 4700
$ws.Range("I17").Value = 4700
$ws.Range("K17").Value = 4700
$ws.Range("M17").Value = -4528
$ws.Range("H100").Value = 851.4167
$ws.Range("I100").Value = 956.9
$ws.Range("K100").Value = 1913.8
$ws.Range("M100").Value = -1372.8
$ws.Range("H107").Value = 3941.5527
$ws.Range("I107").Value = 2751.4583
$ws.Range("K107").Value = 8254.374899999999
$ws.Range("M107").Value = -6334.374899999999
$ws.Range("H132").Value = 21535.611
$ws.Range("I132").Value = 12638.294
$ws.Range("K132").Value = 37914.882
$ws.Range("M132").Value = -35384.882
